$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 16:55:04"
$wsZhCn.Range("H2").Value = "2016-03-19 16:55:20"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 16:55:07"
$wsDeDe.Range("H2").Value = "2016-03-19 16:55:25"
